$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 173
$ws.Range("I12").Value = 173
$ws.Range("K12").Value = 173
$ws.Range("M12").Value = -3

$ws.Range("H33").Value = 220.6842
$ws.Range("I33").Value = 108.21429
$ws.Range("J33").Value = 535.6
$ws.Range("K33").Value = 108.21429
$ws.Range("L33").Value = 535.6
$ws.Range("M33").Value = 120.78571
$ws.Range("N33").Value = -993.6

$ws.Range("H132").Value = 4634172.5
$ws.Range("I132").Value = 5850838.5
$ws.Range("K132").Value = 17552515.5
$ws.Range("M132").Value = -17549985.5

$ws.Range("H137").Value = 1035.2609
$ws.Range("I137").Value = 775.5789
$ws.Range("K137").Value = 2326.7367
$ws.Range("M137").Value = 223.2633000000001

$ws.Range("H138").Value = 1564.0938
$ws.Range("I138").Value = 866.875
$ws.Range("J138").Value = 1627.4773
$ws.Range("K138").Value = 2600.625
$ws.Range("L138").Value = 4882.4319
$ws.Range("M138").Value = 2539.375
$ws.Range("N138").Value = -15162.4319

$ws.Range("H141").Value = 654.625
$ws.Range("I141").Value = 532.34485
$ws.Range("J141").Value = 1836.6666
$ws.Range("K141").Value = 1597.03455
$ws.Range("L141").Value = 5509.9998
$ws.Range("M141").Value = 3582.96545
$ws.Range("N141").Value = -15869.9998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3369.8525
$ws.Range("I32").Value = 3162.228
$ws.Range("K32").Value = 3162.228
$ws.Range("M32").Value = -2875.228

$ws.Range("H45").Value = 1311.6471
$ws.Range("J45").Value = 1413.5
$ws.Range("L45").Value = 1413.5
$ws.Range("N45").Value = -2167.5

$ws.Range("H61").Value = 1216
$ws.Range("I61").Value = 1085.9166
$ws.Range("J61").Value = 1684.3
$ws.Range("K61").Value = 1085.9166
$ws.Range("L61").Value = 1684.3
$ws.Range("M61").Value = -873.9166
$ws.Range("N61").Value = -2108.3

$ws.Range("H74").Value = 1220.5
$ws.Range("I74").Value = 608
$ws.Range("K74").Value = 608
$ws.Range("M74").Value = 266

$ws.Range("H77").Value = 1220.5
$ws.Range("I77").Value = 608
$ws.Range("K77").Value = 3040
$ws.Range("M77").Value = 1328

$ws.Range("H102").Value = 41668092
$ws.Range("I102").Value = 41668092
$ws.Range("K102").Value = 41668092
$ws.Range("M102").Value = -41666470

$ws.Range("H110").Value = 1757.1666
$ws.Range("I110").Value = 1370.4615
$ws.Range("K110").Value = 1370.4615
$ws.Range("M110").Value = 674.5385000000001

$ws.Range("H132").Value = 1641.2174
$ws.Range("I132").Value = 1393.1
$ws.Range("K132").Value = 4179.299999999999
$ws.Range("M132").Value = -1649.299999999999

$ws.Range("H136").Value = 1216
$ws.Range("I136").Value = 1085.9166
$ws.Range("J136").Value = 1684.3
$ws.Range("K136").Value = 3257.7498
$ws.Range("L136").Value = 5052.9
$ws.Range("M136").Value = -707.7498000000001
$ws.Range("N136").Value = -10152.9

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 13158924
$ws.Range("I94").Value = 15625954
$ws.Range("J94").Value = 1433.3334
$ws.Range("K94").Value = 15625954
$ws.Range("L94").Value = 1433.3334
$ws.Range("M94").Value = -15625503
$ws.Range("N94").Value = -2335.3334

$ws.Range("H105").Value = 76925180
$ws.Range("I105").Value = 83335490
$ws.Range("K105").Value = 83335490
$ws.Range("M105").Value = -83333743

$ws.Range("H134").Value = 3768.9583
$ws.Range("I134").Value = 1010.0789
$ws.Range("J134").Value = 14252.7
$ws.Range("K134").Value = 3030.2367
$ws.Range("L134").Value = 42758.10000000001
$ws.Range("M134").Value = -495.2366999999999
$ws.Range("N134").Value = -47828.10000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 83334860
$ws.Range("I16").Value = 100001430
$ws.Range("K16").Value = 100001430
$ws.Range("M16").Value = -100001143

$ws.Range("H31").Value = 1939.9642
$ws.Range("I31").Value = 2005.2174
$ws.Range("J31").Value = 1639.8
$ws.Range("K31").Value = 2005.2174
$ws.Range("L31").Value = 1639.8
$ws.Range("M31").Value = -1710.2174
$ws.Range("N31").Value = -2229.8

$ws.Range("H34").Value = 1939.9642
$ws.Range("I34").Value = 2005.2174
$ws.Range("J34").Value = 1639.8
$ws.Range("K34").Value = 2005.2174
$ws.Range("L34").Value = 1639.8
$ws.Range("M34").Value = -1803.2174
$ws.Range("N34").Value = -2043.8

$ws.Range("H39").Value = 3000
$ws.Range("I39").Value = 3000
$ws.Range("K39").Value = 3000
$ws.Range("M39").Value = -2609

$ws.Range("H49").Value = 3000
$ws.Range("I49").Value = 3000
$ws.Range("K49").Value = 3000
$ws.Range("M49").Value = -2818

$ws.Range("H58").Value = 656.619
$ws.Range("I58").Value = 552.6667
$ws.Range("J58").Value = 1280.3334
$ws.Range("K58").Value = 552.6667
$ws.Range("L58").Value = 1280.3334
$ws.Range("M58").Value = -349.6667
$ws.Range("N58").Value = -1686.3334

$ws.Range("H113").Value = 83334860
$ws.Range("I113").Value = 100001430
$ws.Range("K113").Value = 100001430
$ws.Range("M113").Value = -99999260

$ws.Range("H132").Value = 4813.273
$ws.Range("I132").Value = 5232.346
$ws.Range("J132").Value = 3256.7144
$ws.Range("K132").Value = 15697.038
$ws.Range("L132").Value = 9770.143199999999
$ws.Range("M132").Value = -13167.038
$ws.Range("N132").Value = -14830.1432

$ws.Range("H134").Value = 945.39215
$ws.Range("I134").Value = 941.81396
$ws.Range("K134").Value = 2825.44188
$ws.Range("M134").Value = -290.4418799999999

$ws.Range("H136").Value = 656.619
$ws.Range("I136").Value = 552.6667
$ws.Range("J136").Value = 1280.3334
$ws.Range("K136").Value = 1658.0001
$ws.Range("L136").Value = 3841.0002
$ws.Range("M136").Value = 891.9999
$ws.Range("N136").Value = -8941.0002

$ws.Range("H141").Value = 28960
$ws.Range("J141").Value = 28960
$ws.Range("L141").Value = 28960
$ws.Range("N141").Value = -39320

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1025301.25
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 1025301.25
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 3075903.75
$ws.Range("M4").ClearContents()
$ws.Range("N4").Value = -3076127.75

$ws.Range("H129").Value = 15433588
$ws.Range("I129").Value = 41667188
$ws.Range("J129").Value = 4387861.5
$ws.Range("K129").Value = 125001564
$ws.Range("L129").Value = 13163584.5
$ws.Range("M129").Value = -124996564
$ws.Range("N129").Value = -13173584.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 226.5
$ws.Range("J2").Value = 234.5
$ws.Range("L2").Value = 234.5
$ws.Range("N2").Value = -460.5

$ws.Range("H40").Value = 3509
$ws.Range("J40").Value = 3509
$ws.Range("L40").Value = 3509
$ws.Range("N40").Value = -3811

$ws.Range("H102").Value = 11300
$ws.Range("I102").Value = 9850
$ws.Range("J102").Value = 20000
$ws.Range("K102").Value = 9850
$ws.Range("L102").Value = 20000
$ws.Range("M102").Value = -8228
$ws.Range("N102").Value = -23244

$ws.Range("H107").Value = 707.58826
$ws.Range("I107").Value = 927.6667
$ws.Range("J107").Value = 460
$ws.Range("K107").Value = 927.6667
$ws.Range("L107").Value = 460
$ws.Range("M107").Value = 992.3333
$ws.Range("N107").Value = -4300

$ws.Range("H113").Value = 1723.7693
$ws.Range("I113").Value = 1687.1818
$ws.Range("K113").Value = 1687.1818
$ws.Range("M113").Value = 482.8181999999999

$ws.Range("H122").Value = 2300.5881
$ws.Range("I122").Value = 1442.5834
$ws.Range("J122").Value = 4359.8
$ws.Range("K122").Value = 4327.7502
$ws.Range("L122").Value = 13079.4
$ws.Range("M122").Value = -1877.7502
$ws.Range("N122").Value = -17979.4

$ws.Range("H132").Value = 1733.6333
$ws.Range("I132").Value = 1564.2609
$ws.Range("K132").Value = 4692.7827
$ws.Range("M132").Value = -2162.7827

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H87").Value = 22000
$ws.Range("J87").Value = 22000
$ws.Range("L87").Value = 22000
$ws.Range("N87").Value = -24246

$ws.Range("H90").Value = 22000
$ws.Range("J90").Value = 22000
$ws.Range("L90").Value = 66000
$ws.Range("N90").Value = -77232

$ws.Range("H136").Value = 14663.375
$ws.Range("I136").Value = 21860.4
$ws.Range("J136").Value = 2668.3333
$ws.Range("K136").Value = 65581.20000000001
$ws.Range("L136").Value = 8004.999899999999
$ws.Range("M136").Value = -63031.20000000001
$ws.Range("N136").Value = -13104.9999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").ClearContents()

$ws.Range("H85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").ClearContents()

$ws.Range("H132").Value = 2324.3865
$ws.Range("I132").Value = 2423.389
$ws.Range("K132").Value = 7270.167
$ws.Range("M132").Value = -4740.167

$ws.Range("H136").Value = 606.9524
$ws.Range("I136").Value = 385.66666
$ws.Range("J136").Value = 1160.1666
$ws.Range("K136").Value = 1156.99998
$ws.Range("L136").Value = 3480.4998
$ws.Range("M136").Value = 1393.00002
$ws.Range("N136").Value = -8580.4998
